$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112; this shifts the existing rows 112-133 down to 113-134
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new data record
$ws.Cells.Item(112, 1).Value = 5
$ws.Cells.Item(112, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(112, 3).Value = "Maule"
$ws.Cells.Item(112, 4).Value = 45244
$ws.Cells.Item(112, 5).Value = 7
$ws.Cells.Item(112, 6).Value = 300000000
$ws.Cells.Item(112, 7).Value = "Espárragos"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 3000
$ws.Cells.Item(112, 11).Value = 1400
$ws.Cells.Item(112, 12).Value = 1400
$ws.Cells.Item(112, 13).Value = 1400
$ws.Cells.Item(112, 14).Value = "`$/kilo"
$ws.Cells.Item(112, 15).Value = "Provincia de Linares"
$ws.Cells.Item(112, 16).Value = 1400
$ws.Cells.Item(112, 17).Value = 1
$ws.Cells.Item(112, 18).Value = "Hortaliza"
